$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Accuracy")
$ws.Range("M2").Value = 0.5523012552301255
$ws.Range("M3").Value = 0.6677649643053267
$ws.Range("M4").Value = 0.7246240601503759
$ws.Range("M5").Value = 0.5263157894736842
$ws.Range("M6").Value = 0.5912395492548165
$ws.Range("M7").Value = 0.5906432748538012
$ws.Range("M8").Value = 0.5136986301369864
$ws.Range("M10").Value = 0.6255506607929515
$ws.Range("B11").Value = 0.6997840172786177
$ws.Range("C11").Value = 0.5950323974082073
$ws.Range("D11").Value = 0.5734341252699784
$ws.Range("E11").Value = 0.5604751619870411
$ws.Range("F11").Value = 0.6058315334773218
$ws.Range("G11").Value = 0.642548596112311
$ws.Range("H11").Value = 0.4870410367170626
$ws.Range("I11").Value = 0.6058315334773218
$ws.Range("J11").Value = 0.5043196544276458
$ws.Range("K11").Value = 0.4892008639308855
$ws.Range("L11").Value = 0.6144708423326134
$ws.Range("M11").Value = 0.7505399568034558
$ws.Range("M12").Value = 0.7037383177570093
$ws.Range("M13").Value = 0.6014890282131662

$ws = $wb.Worksheets.Item("Macro_Precision")
$ws.Range("M2").Value = 0.5415139052263731
$ws.Range("M3").Value = 0.6959676126342793
$ws.Range("M4").Value = 0.624845373577437
$ws.Range("M5").Value = 0.5840190700104493
$ws.Range("M6").Value = 0.6350509246071655
$ws.Range("M7").Value = 0.6009005621635501
$ws.Range("M8").Value = 0.5400789715637987
$ws.Range("M9").Value = 0.6304415909802693
$ws.Range("M10").Value = 0.5622785829307568
$ws.Range("B11").Value = 0.7050394725143079
$ws.Range("C11").Value = 0.5716264521894548
$ws.Range("D11").Value = 0.5512988661364021
$ws.Range("E11").Value = 0.559247853676745
$ws.Range("F11").Value = 0.647521959951709
$ws.Range("G11").Value = 0.6453601315002988
$ws.Range("H11").Value = 0.5974310776942355
$ws.Range("I11").Value = 0.5891006128912838
$ws.Range("J11").Value = 0.6031447784810127
$ws.Range("K11").Value = 0.6019638941102756
$ws.Range("L11").Value = 0.6114594692838967
$ws.Range("M11").Value = 0.7455492957746479
$ws.Range("M12").Value = 0.5667899878492364
$ws.Range("M13").Value = 0.4094746437944967

$ws = $wb.Worksheets.Item("Macro_Recall")
$ws.Range("M2").Value = 0.5641493412585743
$ws.Range("M3").Value = 0.6680761826182618
$ws.Range("M4").Value = 0.7642383184834278
$ws.Range("M5").Value = 0.6178977272727273
$ws.Range("M6").Value = 0.6284168195433393
$ws.Range("M7").Value = 0.5941325941325941
$ws.Range("M8").Value = 0.5387295438008679
$ws.Range("M9").Value = 0.6734693877551021
$ws.Range("M10").Value = 0.6151041666666667
$ws.Range("B11").Value = 0.7100525904512616
$ws.Range("C11").Value = 0.54610892982986
$ws.Range("D11").Value = 0.5465020398575215
$ws.Range("E11").Value = 0.560874525990805
$ws.Range("F11").Value = 0.6359058070021526
$ws.Range("G11").Value = 0.649228401720096
$ws.Range("H11").Value = 0.5477101340888717
$ws.Range("I11").Value = 0.5830756545042259
$ws.Range("J11").Value = 0.5600020134903856
$ws.Range("K11").Value = 0.5499297675377409
$ws.Range("L11").Value = 0.5587987132837631
$ws.Range("M11").Value = 0.7507370813018653
$ws.Range("M12").Value = 0.5940479294478528
$ws.Range("M13").Value = 0.4332695954044552

$ws = $wb.Worksheets.Item("Macro_F1")
$ws.Range("M2").Value = 0.5027758323429803
$ws.Range("M3").Value = 0.6556178547749867
$ws.Range("M4").Value = 0.6206690539723331
$ws.Range("M5").Value = 0.5098906100157665
$ws.Range("M6").Value = 0.5903734094235649
$ws.Range("M7").Value = 0.584965847231372
$ws.Range("M8").Value = 0.5130386273970089
$ws.Range("M9").Value = 0.6135211267605634
$ws.Range("M10").Value = 0.5353849117484167
$ws.Range("B11").Value = 0.6989061681263566
$ws.Range("C11").Value = 0.5233977724250098
$ws.Range("D11").Value = 0.5440703719908108
$ws.Range("E11").Value = 0.5569154007761699
$ws.Range("F11").Value = 0.6033665553009047
$ws.Range("G11").Value = 0.64099068233861
$ws.Range("H11").Value = 0.4435060950537389
$ws.Range("I11").Value = 0.5829276629560938
$ws.Range("J11").Value = 0.4720304086157744
$ws.Range("K11").Value = 0.4458492272850917
$ws.Range("L11").Value = 0.5265733414485697
$ws.Range("M11").Value = 0.7467689980075931
$ws.Range("M12").Value = 0.5692436357058248
$ws.Range("M13").Value = 0.4003967869618839
